$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Twitter-scraper related cells
$ws.Range("B2").Value = "#istandwithraeesahkhan"
$ws.Range("E2").Value = "#istandwithraeesahkhan site: twitter.com"
$ws.Range("H2").ClearContents()
$ws.Range("B3").Value = 25
$ws.Range("E3").Value = 75

$ws.Range("I5").Select()
